## AutoCommit_13 июня 2024 г. 10:03:15_SibNout2023
## Re-creates the author's edit: a handful of new "5" marks in the grade
## sheet (J11, J13, E21, K24), the matching conditional-formatting sqref
## growth, and the scrolled/selected cell the author ended the session on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New grade entries.
#    J13 / E21 / K24 need to inherit the same look (green fill + thick
#    side borders) that the sheet already uses for this "total" column,
#    so copy the formatting from an already-styled neighbour before
#    writing the value - that's what Excel's own fill/format-paint does
#    and it keeps the cellXfs entries identical to neighbouring cells.
# ---------------------------------------------------------------------

# J11 -> plain number, no special formatting (matches I11 which is
# also unstyled).
$ws.Range("J11").Value = 5

# J13 -> same look as J17/J28 (the "white" variant used when the I-cell
# in that row is the plain style rather than the green one).
$ws.Range("J17").Copy()
$ws.Range("J13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J13").Value = 5

# E21 -> same look as the rest of row 21's green total cells.
$ws.Range("J12").Copy()
$ws.Range("E21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E21").Value = 5

# K24 -> same look as J24 right next to it.
$ws.Range("J24").Copy()
$ws.Range("K24").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K24").Value = 5

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Conditional formatting: K24 now needs the same color-scale rule
#    that covers the other total cells. The existing rule's AppliesTo
#    already spans many disjoint areas (I4:I12, I29:I31, ... J18); add
#    a matching color-scale rule scoped to K24 so it lights up the same
#    way without disturbing the rule already covering the rest.
#    AddColorScale(3) defaults to Excel's standard red/yellow/green
#    3-color scale, which is exactly what the other cells already use.
# ---------------------------------------------------------------------
$k24 = $ws.Range("K24")
$k24.FormatConditions.AddColorScale(3)

# ---------------------------------------------------------------------
# 3) View state: the author had scrolled the frozen pane down to row 18
#    and left the selection on F21.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 3
$ws.Range("F21").Activate()
